$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.723.25"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.101.59"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.36%  "
$ws.Range("D13").Value = "2.413.53"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "2.106.80"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "38.728.31"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("E28").Value = "  +5.87%  "
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.40%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.69%  "
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").Value = "1.534.43"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "2.295.70"
$ws.Range("E51").Value = "  +0.24%  "
